$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C9 = 33 (this also updates the chart's cached values since Excel
# recalculates series caches from the referenced ranges)
$ws.Range("C9").Value = 33

# Move the active-cell selection from C12 to C11
$ws.Range("C11").Select()
